$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells that look like plain numbers stay as text,
# matching the source data which stores all prices as literal strings.
$ws.Range('D2').Value = '34.059.21'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.789.82'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.98'
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.545'
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.29'
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  +3.92%  '
$ws.Range('E10').Value = '  -3.21%  '
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('D12').Value = '2.047.00'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.46'
$ws.Range('E13').Value = '  +4.55%  '
$ws.Range('D14').Value = '1.791.28'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.624'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').Value = '34.050.94'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.06'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '242.51'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').Value = '0.0₃0775'
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.77'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('E24').Value = '  -3.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.98'
$ws.Range('E25').Value = '  +2.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.16'
$ws.Range('E26').Value = '  +1.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.22'
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('E28').Value = '  +1.20%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  +2.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0518'
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('E32').Value = '  -0.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.64'
$ws.Range('E33').Value = '  +4.47%  '
$ws.Range('E34').Value = '  +1.77%  '
$ws.Range('D35').Value = '1.403.22'
$ws.Range('E35').Value = '  +1.44%  '
$ws.Range('E36').Value = '  +1.21%  '
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0189'
$ws.Range('E38').Value = '  +2.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.34'
$ws.Range('E39').Value = '  +7.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '80.02'
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.920'
$ws.Range('E42').Value = '  +0.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.72'
$ws.Range('E43').Value = '  +14.79%  '
$ws.Range('E44').Value = '  -0.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.13'
$ws.Range('E45').Value = '  +4.53%  '
$ws.Range('E46').Value = '  +3.60%  '
$ws.Range('E47').Value = '  +1.92%  '
$ws.Range('E48').Value = '  +2.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '107.54'
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('D50').Value = '1.948.83'
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('E51').Value = '  -0.06%  '
